# Corrected spelling mistakes in ppt (Module 6.1 - Control logic design)
#
# Helper: replace the text of a TextRange sub-range with brand-new text while
# collapsing the result down to a single run that carries clean (non-"err",
# non-split) run properties. PowerPoint's TextRange.Text setter keeps the old
# run boundaries/formatting when the new string shares a common prefix/suffix
# with the old one (so misspelled words flagged with err="1" stick around and
# runs stay fragmented). Writing an unrelated placeholder first breaks that
# prefix/suffix match, so the following real assignment collapses cleanly to
# one run using the first old run's properties.
function Set-RangeText($range, $newText) {
    $range.Text = "~~tmp~~"
    $range.Text = $newText
}

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 13 - "Design of hard wired control"
# ---------------------------------------------------------------------
$s13 = $p.Slides.Item(13)
$sh13 = $s13.Shapes.Item(1)
$tr13 = $sh13.TextFrame.TextRange

Set-RangeText ($tr13.Paragraphs(3, 1)) "So 8 D flip-flops used with the following input functions (implemented in decision logic). "
Set-RangeText ($tr13.Paragraphs(4, 1)) "At a time, any of the state will be active and corresponding D flipflop will be activated. "
Set-RangeText ($tr13.Paragraphs(5, 1)) "The output signals are generated according to the Boolean function (in table) using OR gate."

# ---------------------------------------------------------------------
# Slide 28 - "G1 and G2 -flipflops(sequence reg)" / "Input to PLA- ... input"
# ---------------------------------------------------------------------
$s28 = $p.Slides.Item(28)
$sh28 = $s28.Shapes.Item(2)
$tr28 = $sh28.TextFrame.TextRange

# "-flipflops(sequence " (3 runs) -> "-flip-flops (sequence " (1 run)
$para6 = $tr28.Paragraphs(6, 1)
Set-RangeText ($para6.Characters(11, 20)) "–flip-flops (sequence "

# " input" -> " inputs"
$para7 = $tr28.Paragraphs(7, 1)
Set-RangeText ($para7.Characters(47, 6)) " inputs"
